$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data (and two ranking swaps)

# Row 2
$ws.Cells.Item(2, 4).Value = '59.342.93'
$ws.Cells.Item(2, 5).Value = '  +1.41%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.326.77'
$ws.Cells.Item(3, 5).Value = '  +0.98%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''1.00'
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 5).Value = '  -0.08%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''543.43'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  -0.70%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''132.46'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  +0.48%  '

# Row 7
$ws.Cells.Item(7, 4).Value = '''1.00'
$ws.Cells.Item(7, 4).ClearFormats()
$ws.Cells.Item(7, 5).Value = '  -0.04%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  +2.15%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '2.319.51'
$ws.Cells.Item(9, 5).Value = '  +0.77%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  -0.47%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.25%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''0.151'
$ws.Cells.Item(12, 4).ClearFormats()
$ws.Cells.Item(12, 5).Value = '  +0.62%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''0.332'
$ws.Cells.Item(13, 4).ClearFormats()
$ws.Cells.Item(13, 5).Value = '  +0.27%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +0.88%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '2.736.24'
$ws.Cells.Item(15, 5).Value = '  +0.81%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '59.236.60'
$ws.Cells.Item(16, 5).Value = '  +1.24%  '

# Row 17
$ws.Cells.Item(17, 5).Value = '  +0.32%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '2.327.48'
$ws.Cells.Item(18, 5).Value = '  +2.48%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  +0.31%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '''4.19'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  -1.91%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''314.54'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +0.46%  '

# Row 22
$ws.Cells.Item(22, 4).Value = '''6.64'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +3.07%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.08%  '

# Row 24
$ws.Cells.Item(24, 5).Value = '  -0.49%  '

# Row 25
$ws.Cells.Item(25, 4).Value = '''0.174'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +3.55%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''1.00'
$ws.Cells.Item(26, 4).ClearFormats()
$ws.Cells.Item(26, 5).Value = '  +0.12%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''8.01'
$ws.Cells.Item(27, 4).ClearFormats()
$ws.Cells.Item(27, 5).Value = '  +0.45%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''1.32'
$ws.Cells.Item(28, 4).ClearFormats()
$ws.Cells.Item(28, 5).Value = '  +1.56%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Monero'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(29, 4).Value = '''171.40'
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +0.41%  '

# Row 30
$ws.Cells.Item(30, 2).Value = 'PancakeSwap'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(30, 4).Value = '''1.73'
$ws.Cells.Item(30, 4).ClearFormats()
$ws.Cells.Item(30, 5).Value = '  -1.98%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +8.79%  '

# Row 32
$ws.Cells.Item(32, 4).Value = '0.0₃0741'
$ws.Cells.Item(32, 5).Value = '  +3.15%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''5.89'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  +2.86%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'ImmutableX'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(34, 4).Value = '''1.44'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  +16.30%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'PolygonEcosystemToken'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Cells.Item(35, 4).Value = '''0.385'
$ws.Cells.Item(35, 4).ClearFormats()
$ws.Cells.Item(35, 5).Value = '  +1.90%  '

# Row 36
$ws.Cells.Item(36, 5).Value = '  +0.03%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''17.89'
$ws.Cells.Item(37, 4).ClearFormats()
$ws.Cells.Item(37, 5).Value = '  +0.80%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -0.09%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  +3.68%  '

# Row 40
$ws.Cells.Item(40, 4).Value = '''318.40'
$ws.Cells.Item(40, 4).ClearFormats()
$ws.Cells.Item(40, 5).Value = '  +9.12%  '

# Row 41
$ws.Cells.Item(41, 2).Value = 'Stacks'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(41, 4).Value = '''1.53'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +2.56%  '

# Row 42
$ws.Cells.Item(42, 2).Value = 'OKB'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(42, 4).Value = '''38.01'
$ws.Cells.Item(42, 4).ClearFormats()
$ws.Cells.Item(42, 5).Value = '  -0.28%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''142.68'
$ws.Cells.Item(43, 4).ClearFormats()
$ws.Cells.Item(43, 5).Value = '  +2.02%  '

# Row 44
$ws.Cells.Item(44, 4).Value = '''3.45'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = '  +1.04%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  +0.63%  '

# Row 46
$ws.Cells.Item(46, 4).Value = '''0.0494'
$ws.Cells.Item(46, 4).ClearFormats()
$ws.Cells.Item(46, 5).Value = '  -1.17%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '''0.558'
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = '  -0.86%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''18.43'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = '  +1.16%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -1.08%  '

# Row 50
$ws.Cells.Item(50, 5).Value = '  +0.02%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -0.30%  '
